$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 74 values that changed
$ws.Range("J74").Value = -570
$ws.Range("K74").Value = 71
$ws.Range("L74").Value = -1044
$ws.Range("M74").Value = -152
$ws.Range("N74").Value = 650
$ws.Range("O74").Value = 80
$ws.Range("S74").Value = -1410
$ws.Range("V74").Value = 3
$ws.Range("W74").Value = -158
$ws.Range("X74").Value = -4
$ws.Range("Z74").Value = 363
$ws.Range("AA74").Value = 18

# Add new row 75 with new quarter data
# (force text storage for the date-like label so it matches the
#  existing "dd-mm-yyyy" text entries in column A instead of being
#  auto-converted to a date serial number)
$ws.Range("A75").NumberFormat = "@"
$ws.Range("A75").Value = "01-04-2021"
$ws.Range("A75").ClearFormats()
$ws.Range("B75").Value = -93
$ws.Range("C75").Value = -93
$ws.Range("D75").Value = 1
$ws.Range("E75").Value = -1
$ws.Range("F75").Value = 9016
$ws.Range("G75").Value = -78
$ws.Range("H75").Value = 267
$ws.Range("I75").Value = -417
$ws.Range("J75").Value = 851
$ws.Range("K75").Value = -91
$ws.Range("L75").Value = -443
$ws.Range("M75").Value = 727
$ws.Range("N75").Value = -652
$ws.Range("O75").Value = 8865
$ws.Range("P75").Value = -12
$ws.Range("Q75").Value = 0
$ws.Range("R75").Value = 0
$ws.Range("S75").Value = 3882
$ws.Range("T75").Value = 1533
$ws.Range("U75").Value = 1300
$ws.Range("V75").Value = 15
$ws.Range("W75").Value = -2
$ws.Range("X75").Value = 77
$ws.Range("Y75").Value = 85
$ws.Range("Z75").Value = 815
$ws.Range("AA75").Value = 59
